$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: "年龄" -> "年齡"
$ws.Range("B1").Value = "年齡"

# Row 2: "狗子"/"18" -> "小明"/"11"
$ws.Range("A2").Value = "小明"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "11"

# New row 3: "小华"/"12"
$ws.Range("A3").Value = "小华"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "12"
